# Applies the cryptos-list refresh described by the commit diff.
# Numeric-looking Price (column D) values are given a leading apostrophe so
# Excel keeps them as literal text (matching the source inlineStr cells)
# instead of auto-converting them to numbers and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.023.36"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.236.57"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'249.87"
$ws.Range("E5").Value = "  +7.02%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'72.12"
$ws.Range("E7").Value = "  +4.23%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "  +6.18%  "
$ws.Range("D10").Value = "'41.41"
$ws.Range("E10").Value = "  +16.95%  "
$ws.Range("D11").Value = "'0.0980"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "'58.10"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  +6.94%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "2.571.34"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "'15.04"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "'0.869"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "2.231.93"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "42.021.84"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").Value = "0.0₃0977"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "'6.27"
$ws.Range("D22").Value = "'73.34"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'236.67"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  +10.21%  "
$ws.Range("D25").Value = "'3.98"
$ws.Range("E25").Value = "  +8.86%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'2.54"
$ws.Range("E27").Value = "  +8.20%  "
$ws.Range("D28").Value = "'10.80"
$ws.Range("E28").Value = "  +7.56%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'171.71"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").Value = "'20.92"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").Value = "'0.125"
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").Value = "'5.55"
$ws.Range("E34").Value = "  +4.74%  "
$ws.Range("D35").Value = "'0.0734"
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("D36").Value = "'4.75"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").Value = "'26.48"
$ws.Range("E37").Value = "  +24.78%  "
$ws.Range("D38").Value = "'4.03"
$ws.Range("E38").Value = "  +11.67%  "
$ws.Range("D39").Value = "'0.0300"
$ws.Range("E39").Value = "  +12.83%  "
$ws.Range("D40").Value = "'2.31"
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").Value = "'6.04"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("D42").Value = "'68.27"
$ws.Range("E42").Value = "  +3.32%  "
$ws.Range("D43").Value = "'12.25"
$ws.Range("E43").Value = "  +23.46%  "
$ws.Range("E44").Value = "  +10.89%  "
$ws.Range("D45").Value = "'4.93"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").Value = "'8.82"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "'0.102"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").Value = "'4.70"
$ws.Range("E48").Value = "  +6.87%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  +8.27%  "
$ws.Range("E51").Value = "  +1.92%  "
